# Add data for 2022-05-09 (carjacking by neighborhood by month)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-05-09"

# Update the column header label for the current (partial) month
$ws.Range("B1").Value = "May 2022 (through May 09)"

# --- Cell value updates / additions (row = neighborhood, column = month) ---

# Row 2 - Englewood
$ws.Range("L2").Value = 4          # May 2020: 3 -> 4

# Row 3 - Austin
$ws.Range("L3").Value = 1          # May 2020: new
$ws.Range("AK3").Value = 2         # May 2015: new

# Row 4 - Humboldt Park
$ws.Range("AA4").Value = 1         # May 2017: new

# Row 11 - Roseland
$ws.Range("B11").Value = 3         # May 2022: 2 -> 3

# Row 12 - Kenwood
$ws.Range("G12").Value = 2         # May 2021: 1 -> 2

# Row 13 - Washington Heights
$ws.Range("V13").Value = 1         # May 2018: new

# Row 15 - Lake View
$ws.Range("B15").Value = 2         # May 2022: 1 -> 2

# Row 22 - Brighton Park
$ws.Range("B22").Value = 4         # May 2022: 3 -> 4

# Row 27 - Wicker Park
$ws.Range("L27").Value = 1         # May 2020: new

# Row 35 - Avondale
$ws.Range("Q35").Value = 1         # May 2019: new

# Row 39 - New City
$ws.Range("G39").Value = 2         # May 2021: 1 -> 2

# Row 46 - Little Village
$ws.Range("G46").Value = 2         # May 2021: new

# Row 57 - Chinatown
$ws.Range("AA57").Value = 2        # May 2017: 1 -> 2

# Row 76 - North Center
$ws.Range("AA76").Value = 1        # May 2017: new
